$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '30.265.08'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.868.64'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '244.62'
$ws.Range("E5").Value = '  +4.26%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.4721'
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.2873'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '42.66'
$ws.Range("E9").Value = '  -2.91%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.06469'
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '21.07'
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.929.73'
$ws.Range("E12").Value = '  +3.31%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.07754'
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '94.91'
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.7093'
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '5.098'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '275.61'
$ws.Range("E17").Value = '  +2.69%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '30.266.93'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '13.35'
$ws.Range("E19").Value = '  -4.44%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '0.000007557'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '2.128.30'
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '5.212'
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '6.138'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '9.265'
$ws.Range("E26").Value = '  -1.52%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '165.12'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '18.87'
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '1.904'
$ws.Range("E29").Value = '  -2.21%  '
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.09850'
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '1.515'
$ws.Range("E32").Value = '  +3.66%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '4.250'
$ws.Range("E33").Value = '  -3.00%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '4.032'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.04756'
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '1.122'
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.6932'
$ws.Range("E37").Value = '  -1.42%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '2.708'
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.01842'
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '2.744'
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '6.311'
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '70.38'
$ws.Range("E42").Value = '  -4.23%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.8419'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '1.903'
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.4102'
$ws.Range("E46").Value = '  -1.82%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '101.95'
$ws.Range("E47").Value = '  -1.31%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '9.276'
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '7.089'
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '35.21'
$ws.Range("E50").Value = '  +1.98%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '917.92'
$ws.Range("E51").Value = '  -6.60%  '
